$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- 1. Remove the "Meta description: ..." paragraph that sits right after
#        the Heading1 title at the top of the document. ---
$metaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Meta description:*") {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# --- 2. Before the final paragraph (the former "feature image prompt"
#        paragraph), insert a new bold paragraph repeating the page title. ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphBefore()
$titlePara = $d.Paragraphs($d.Paragraphs.Count - 1)
$null = $titlePara.Range.InsertXML("<w:p $wNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Asgardian Stones for Free - Review and Features</w:t></w:r></w:p>")

# --- 3. Replace the text of the final paragraph with the meta description
#        copy, keeping its italic formatting. ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$null = $lastPara.Range.InsertXML("<w:p $wNs><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of the Asgardian Stones slot game featuring gameplay, bonus features, graphics and sound, and betting options. Play for free today.</w:t></w:r></w:p>")
